$d = $word.ActiveDocument

$pairs = @(
    @("149×8=", "954×8="),
    @("358×3=", "371×4="),
    @("653×2=", "354×6="),
    @("411×2=", "370×8="),
    @("344×7=", "912×3="),
    @("622×9=", "591×8="),
    @("115×3=", "332×6="),
    @("289×6=", "396×8="),
    @("162×4=", "868×8="),
    @("931×3=", "825×7="),
    @("297×6=", "968×7="),
    @("331×7=", "762×8="),
    @("674×4=", "307×6="),
    @("854×6=", "158×4="),
    @("542×3=", "420×3="),
    @("174×8=", "132×2="),
    @("154×4=", "247×9="),
    @("655×3=", "686×8="),
    @("911×4=", "890×6="),
    @("831×9=", "783×8="),
    @("788×9=", "586×7="),
    @("314×3=", "759×2="),
    @("794×3=", "467×9="),
    @("792×4=", "198×3="),
    @("629×3=", "548×7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
